$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.942.63"
$ws.Range("E2").Value = "  +3.16%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.572.02"
$ws.Range("E3").Value = "  +0.50%  "
$ws.Range("E4").Value = "  -1.28%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.71"
$ws.Range("E5").Value = "  +0.57%  "
$ws.Range("E6").Value = "  +0.27%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.995"
$ws.Range("E7").Value = "  -1.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.12"
$ws.Range("E8").Value = "  +5.39%  "
$ws.Range("E9").Value = "  +0.35%  "
$ws.Range("E10").Value = "  -0.36%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0879"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.795.98"
$ws.Range("E12").Value = "  +0.42%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.566.13"
$ws.Range("E13").Value = "  +0.19%  "
$ws.Range("E14").Value = "  -0.57%  "
$ws.Range("E15").Value = "  +0.28%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.905.71"
$ws.Range("E16").Value = "  +3.04%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.33"
$ws.Range("E17").Value = "  +2.16%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "228.10"
$ws.Range("E18").Value = "  +6.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0702"
$ws.Range("E19").Value = "  +0.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.43"
$ws.Range("E20").Value = "  +0.88%  "
$ws.Range("E21").Value = "  -1.23%  "
$ws.Range("E22").Value = "  -0.89%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.27"
$ws.Range("E23").Value = "  +0.83%  "
$ws.Range("E24").Value = "  -0.32%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.51"
$ws.Range("E25").Value = "  -1.50%  "
$ws.Range("E26").Value = "  +0.83%  "
$ws.Range("E27").Value = "  -0.57%  "
$ws.Range("E28").Value = "  +0.25%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.996"
$ws.Range("E29").Value = "  -1.10%  "
$ws.Range("E30").Value = "  +0.22%  "
$ws.Range("E31").Value = "  +0.25%  "
$ws.Range("E32").Value = "  -0.57%  "
$ws.Range("E33").Value = "  -1.64%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.409.37"
$ws.Range("E34").Value = "  -2.14%  "
$ws.Range("E35").Value = "  -1.42%  "
$ws.Range("E36").Value = "  -3.79%  "
$ws.Range("E37").Value = "  -1.88%  "
$ws.Range("E38").Value = "  -0.26%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.537"
$ws.Range("E39").Value = "  +1.07%  "
$ws.Range("E40").Value = "  +3.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.804"
$ws.Range("E41").Value = "  -0.13%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.994"
$ws.Range("E42").Value = "  -1.24%  "
$ws.Range("E43").Value = "  -3.51%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.972"
$ws.Range("E44").Value = "  -2.60%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.80"
$ws.Range("E45").Value = "  +3.68%  "
$ws.Range("E46").Value = "  -1.57%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.706.84"
$ws.Range("E47").Value = "  +0.43%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.56"
$ws.Range("E48").Value = "  +0.70%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₆0105"
$ws.Range("E49").Value = "  +2.20%  "
$ws.Range("E50").Value = "  +1.13%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0939"
$ws.Range("E51").Value = "  -1.85%  "
